# AfDD 2022 Annex Table 14 (Digitalisation) refresh:
#  - Header labels: survey reference window moves from 2011-20 to 2012-21.
#  - CAF (row 16) and Rwanda (row 31) enterprise-survey shares become unavailable ('..').
#  - Downstream country/region aggregates in C:J are recalculated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab14")

# --- Row 2 column headers: update the reference period in each label ---
$ws.Range("C2").Value = "Percentage of large firms using their own website, most recent measure 2012-21"
$ws.Range("D2").Value = "Percentage of medium-sized firms using their own website, most recent measure 2012-21"
$ws.Range("E2").Value = "Percentage of small firms using their own website, most recent measure 2012-21"
$ws.Range("F2").Value = "Percentage of firms of all sizes using their own website, most recent measure 2012-21"
$ws.Range("G2").Value = "Percentage of large firms using email to interact with clients/suppliers, most recent measure 2012-21"
$ws.Range("H2").Value = "Percentage of medium-sized firms using email to interact with clients/suppliers, most recent measure 2012-21"
$ws.Range("I2").Value = "Percentage of small firms using email to interact with clients/suppliers, most recent measure 2012-21"
$ws.Range("J2").Value = "Percentage of firms of all sizes using email to interact with clients/suppliers, most recent measure 2012-21"

# --- Row 16 (Central African Republic): website/e-mail shares withdrawn -> '..' ---
$ws.Range("C16").Value = ".."
$ws.Range("D16").Value = ".."
$ws.Range("E16").Value = ".."
$ws.Range("F16").Value = ".."
$ws.Range("G16").Value = ".."
$ws.Range("H16").Value = ".."
$ws.Range("I16").Value = ".."
$ws.Range("J16").Value = ".."

# --- Row 31 (Rwanda): website/e-mail shares withdrawn -> '..' ---
$ws.Range("C31").Value = ".."
$ws.Range("D31").Value = ".."
$ws.Range("E31").Value = ".."
$ws.Range("F31").Value = ".."
$ws.Range("G31").Value = ".."
$ws.Range("H31").Value = ".."
$ws.Range("I31").Value = ".."
$ws.Range("J31").Value = ".."

# --- Recalculated aggregate figures (columns C:J) ---
# Row 23: Central Africa
$ws.Range("C23").Value = 56.275
$ws.Range("D23").Value = 28.875
$ws.Range("E23").Value = 13.8
$ws.Range("F23").Value = 18.675
$ws.Range("G23").Value = 85.8
$ws.Range("H23").Value = 68.2
$ws.Range("I23").Value = 36.35
$ws.Range("J23").Value = 44.425
# Row 38: East Africa
$ws.Range("C38").Value = 71.4375
$ws.Range("D38").Value = 52.0125
$ws.Range("E38").Value = 27.35
$ws.Range("F38").Value = 35.55
$ws.Range("G38").Value = 85.6375
$ws.Range("H38").Value = 74.7375
$ws.Range("I38").Value = 50.4625
$ws.Range("J38").Value = 58.5
# Row 62: Africa
$ws.Range("C62").Value = 69.1
$ws.Range("D62").Value = 46.9142857142857
$ws.Range("E62").Value = 24.3285714285714
$ws.Range("F62").Value = 33.18
$ws.Range("G62").Value = 89.0848484848485
$ws.Range("H62").Value = 76.3848484848485
$ws.Range("I62").Value = 52.1818181818182
$ws.Range("J62").Value = 60.569696969697
# Row 63: World outside Africa
$ws.Range("C63").Value = 80.3821428571429
$ws.Range("D63").Value = 60.7535714285714
$ws.Range("E63").Value = 39.3571428571429
$ws.Range("F63").Value = 48.5892857142857
$ws.Range("G63").Value = 93.8777777777778
$ws.Range("H63").Value = 83.1555555555556
$ws.Range("I63").Value = 64.9592592592593
$ws.Range("J63").Value = 72.5111111111111
# Row 65: Asia (no high inc.)
$ws.Range("C65").Value = 75.6666666666667
$ws.Range("D65").Value = 47.325
$ws.Range("E65").Value = 27.3583333333333
$ws.Range("F65").Value = 36.925
$ws.Range("G65").Value = 92.875
$ws.Range("H65").Value = 70.3916666666667
$ws.Range("I65").Value = 48.7833333333333
$ws.Range("J65").Value = 58.55
# Row 66: World
$ws.Range("C66").Value = 74.1142857142857
$ws.Range("D66").Value = 53.0650793650794
$ws.Range("E66").Value = 31.0079365079365
$ws.Range("F66").Value = 40.0285714285714
$ws.Range("G66").Value = 91.2416666666667
$ws.Range("H66").Value = 79.4316666666667
$ws.Range("I66").Value = 57.9316666666667
$ws.Range("J66").Value = 65.9433333333333
# Row 67: COMESA
$ws.Range("C67").Value = 71.7285714285714
$ws.Range("D67").Value = 51.9714285714286
$ws.Range("E67").Value = 32.1785714285714
$ws.Range("F67").Value = 40.2285714285714
$ws.Range("G67").Value = 90.3142857142857
$ws.Range("H67").Value = 79.1857142857143
$ws.Range("I67").Value = 57.1071428571429
$ws.Range("J67").Value = 65.4071428571429
# Row 68: CEN-SAD
$ws.Range("C68").Value = 72.9888888888889
$ws.Range("D68").Value = 48.6666666666667
$ws.Range("E68").Value = 25.4444444444445
$ws.Range("F68").Value = 35.7111111111111
$ws.Range("G68").Value = 91.1
$ws.Range("H68").Value = 77.4529411764706
$ws.Range("I68").Value = 53.9411764705882
$ws.Range("J68").Value = 62.6882352941177
# Row 69: EAC
$ws.Range("C69").Value = 69.8
$ws.Range("D69").Value = 46.58
$ws.Range("E69").Value = 20.98
$ws.Range("F69").Value = 28.06
$ws.Range("G69").Value = 85.96
$ws.Range("H69").Value = 70.08
$ws.Range("I69").Value = 43.94
$ws.Range("J69").Value = 52.08
# Row 70: ECCAS
$ws.Range("C70").Value = 56.275
$ws.Range("D70").Value = 28.875
$ws.Range("E70").Value = 13.8
$ws.Range("F70").Value = 18.675
$ws.Range("G70").Value = 85.8
$ws.Range("H70").Value = 68.2
$ws.Range("I70").Value = 36.35
$ws.Range("J70").Value = 44.425
# Row 81: ROW, Resource-rich countries
$ws.Range("C81").Value = 78.2166666666667
$ws.Range("D81").Value = 57.85
$ws.Range("E81").Value = 40.3166666666667
$ws.Range("F81").Value = 46.5833333333333
$ws.Range("G81").Value = 91.84
$ws.Range("H81").Value = 75.12
$ws.Range("I81").Value = 57.86
$ws.Range("J81").Value = 63.3
# Row 82: Africa, Non-resource-rich countries
$ws.Range("C82").Value = 69.296875
$ws.Range("D82").Value = 47.771875
$ws.Range("E82").Value = 25.184375
$ws.Range("F82").Value = 34.028125
$ws.Range("G82").Value = 89.4166666666667
$ws.Range("H82").Value = 77.9566666666667
$ws.Range("I82").Value = 53.7733333333334
$ws.Range("J82").Value = 62.27
# Row 83: ROW, Non-resource-rich countries
$ws.Range("C83").Value = 80.9727272727273
$ws.Range("D83").Value = 61.5454545454545
$ws.Range("E83").Value = 39.0954545454546
$ws.Range("F83").Value = 49.1363636363636
$ws.Range("G83").Value = 94.3409090909091
$ws.Range("H83").Value = 84.9818181818182
$ws.Range("I83").Value = 66.5727272727273
$ws.Range("J83").Value = 74.6045454545455
# Row 84: Africa, Low income countries
$ws.Range("C84").Value = 63.6058823529412
$ws.Range("D84").Value = 41.2411764705882
$ws.Range("E84").Value = 20.3823529411765
$ws.Range("F84").Value = 28.1058823529412
$ws.Range("G84").Value = 85.6533333333334
$ws.Range("H84").Value = 73.5466666666667
$ws.Range("I84").Value = 47.7066666666667
$ws.Range("J84").Value = 56.0266666666667
# Row 87: ROW, Lower middle income countries
$ws.Range("C87").Value = 77.525
$ws.Range("D87").Value = 51.4333333333333
$ws.Range("E87").Value = 28.2416666666667
$ws.Range("F87").Value = 38.8833333333333
$ws.Range("G87").Value = 93.5416666666667
$ws.Range("H87").Value = 76.825
$ws.Range("I87").Value = 55.3
$ws.Range("J87").Value = 65.2083333333334
# Row 89: ROW, Upper middle income countries
$ws.Range("C89").Value = 84.1357142857143
$ws.Range("D89").Value = 71.2857142857143
$ws.Range("E89").Value = 52.4142857142857
$ws.Range("F89").Value = 60.7928571428572
$ws.Range("G89").Value = 93.6230769230769
$ws.Range("H89").Value = 91.2307692307692
$ws.Range("I89").Value = 78.4384615384615
$ws.Range("J89").Value = 83.7384615384616
# Row 91: Africa, Least Developed Countries
$ws.Range("C91").Value = 65.2666666666667
$ws.Range("D91").Value = 41.9291666666667
$ws.Range("E91").Value = 19.4208333333333
$ws.Range("F91").Value = 28.7416666666667
$ws.Range("G91").Value = 86.5590909090909
$ws.Range("H91").Value = 73.7636363636364
$ws.Range("I91").Value = 48.3636363636364
$ws.Range("J91").Value = 57.1318181818182
# Row 95: Africa, Landlocked Developing Countries
$ws.Range("C95").Value = 66.025
$ws.Range("D95").Value = 41.8166666666667
$ws.Range("E95").Value = 24.2833333333333
$ws.Range("F95").Value = 30.9333333333333
$ws.Range("G95").Value = 85.7583333333333
$ws.Range("H95").Value = 75.0333333333334
$ws.Range("I95").Value = 51.3416666666667
$ws.Range("J95").Value = 59.0333333333334
# Row 97: Africa, Fragile States
$ws.Range("C97").Value = 64.6740740740741
$ws.Range("D97").Value = 42.1444444444445
$ws.Range("E97").Value = 21.7111111111111
$ws.Range("F97").Value = 29.0185185185185
$ws.Range("G97").Value = 86.856
$ws.Range("H97").Value = 72.688
$ws.Range("I97").Value = 47.508
$ws.Range("J97").Value = 55.788
# Row 98: ROW, Fragile States
$ws.Range("C98").Value = 77.525
$ws.Range("D98").Value = 50.25
$ws.Range("E98").Value = 23.95
$ws.Range("F98").Value = 35.725
$ws.Range("G98").Value = 94.7625
$ws.Range("H98").Value = 71.775
$ws.Range("I98").Value = 48.5875
$ws.Range("J98").Value = 59.6
# Row 99: Extremely Fragile States
$ws.Range("C99").Value = 65.2428571428572
$ws.Range("D99").Value = 41.8857142857143
$ws.Range("E99").Value = 21.2142857142857
$ws.Range("F99").Value = 26.4571428571429
$ws.Range("G99").Value = 85.3428571428572
$ws.Range("H99").Value = 70.4
$ws.Range("I99").Value = 39.6428571428572
$ws.Range("J99").Value = 46.8857142857143
